$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new location ("JFK Parkway", CO) was added to the data table. It was
# inserted as a new row 10, pushing the existing rows 10-12 (and the blank
# spacer rows below) down by one.
$ws.Rows("10:10").Insert()

$ws.Range("A10").Value = "JFK Parkway"
$ws.Range("C10").Value = "CO"
$ws.Range("D10").Value = 40.525444999999998
$ws.Range("E10").Value = -105.07343899999999

# Update the active selection to reflect where the editor left off.
$ws.Range("D22").Select() | Out-Null
